$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AML-CFT")

# The country / central-bank reference table that used to live in
# B20:C30 was removed (values only — the cell styles/borders for
# those rows stay as they were).
$ws.Range("B20:C30").ClearContents()

# With the text gone, the rows that used to wrap (and so had an
# explicit 30pt height) collapse back down to the default height.
$ws.Rows("20:30").AutoFit()

# Bring the AML-CFT tab to the front and leave the view scrolled /
# selected over the area that used to hold the table.
$ws.Activate()
$ws.Range("A20:D31").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
